# Updated cryptos list with latest price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.898.80"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").Value = "1.949.20"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.87"
$ws.Range("E5").Value = "  -2.56%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4874"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2924"
$ws.Range("E8").Value = "  -1.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06852"
$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  +0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "106.08"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "1.948.85"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07740"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.307"
$ws.Range("E14").Value = "  -2.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6933"
$ws.Range("E15").Value = "  -3.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.90"
$ws.Range("E16").Value = "  -2.75%  "

$ws.Range("D17").Value = "30.903.14"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.14"
$ws.Range("E18").Value = "  -1.58%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007696"
$ws.Range("E19").Value = "  -1.04%  "

$ws.Range("D20").Value = "2.208.01"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.428"
$ws.Range("E22").Value = "  -4.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.460"
$ws.Range("E24").Value = "  -2.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.682"
$ws.Range("E25").Value = "  -3.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.43"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.61"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.156"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1035"
$ws.Range("E29").Value = "  -3.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  -4.02%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.553"
$ws.Range("E31").Value = "  -3.07%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.546"
$ws.Range("E32").Value = "  -5.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.357"
$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04839"
$ws.Range("E34").Value = "  -4.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7457"
$ws.Range("E35").Value = "  -3.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").Value = "  -1.66%  "

$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01986"
$ws.Range("E38").Value = "  -3.42%  "

$ws.Range("E39").Value = "  -1.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.475"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.55"
$ws.Range("E41").Value = "  +3.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.082"
$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8970"

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.91"
$ws.Range("E44").Value = "  -1.87%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4401"
$ws.Range("E45").Value = "  -1.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9992"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.677"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "991.45"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1242"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.199"
$ws.Range("E50").Value = "  -2.80%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.62"
$ws.Range("E51").Value = "  -1.53%  "
